$wb = $excel.ActiveWorkbook

# The same update applies to both the "展览" and "全部类型" worksheets,
# which contain duplicated data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2273
    $ws.Range("F3").Value = 1721
}
